# Auto-generated edit script: update cryptos list data (price & volume) per commit
# 'Updated cryptos list on Wed Aug 28 18:34:00 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.848.99'
$ws.Range('E2').Value = '  -4.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.520.84'
$ws.Range('E3').Value = '  -2.05%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '533.56'
$ws.Range('E5').Value = '  -2.90%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.48'
$ws.Range('E6').Value = '  -7.41%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -4.43%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.517.47'
$ws.Range('E9').Value = '  -2.34%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0996'
$ws.Range('E10').Value = '  -4.34%  '
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.58'
$ws.Range('E12').Value = '  +2.45%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.350'
$ws.Range('E13').Value = '  -4.49%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.959.53'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.79'
$ws.Range('E15').Value = '  -7.14%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '58.828.73'
$ws.Range('E16').Value = '  -4.79%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000138'
$ws.Range('E17').Value = '  -4.72%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.511.61'
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  -2.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.27'
$ws.Range('E20').Value = '  -6.32%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '321.57'
$ws.Range('E21').Value = '  -4.62%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.74'
$ws.Range('E23').Value = '  -5.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.69'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.437'
$ws.Range('E25').Value = '  -11.36%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.162'
$ws.Range('E26').Value = '  -3.58%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.993'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.607.17'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.72'
$ws.Range('E29').Value = '  -5.33%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.86'
$ws.Range('E30').Value = '  -7.57%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0768'
$ws.Range('E31').Value = '  -8.55%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.25'
$ws.Range('E32').Value = '  -6.75%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.78'
$ws.Range('E33').Value = '  -6.56%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '156.58'
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.40'
$ws.Range('E36').Value = '  -2.03%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.47'
$ws.Range('E37').Value = '  -3.88%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.37'
$ws.Range('E38').Value = '  -9.55%  '
$ws.Range('E39').Value = '  -11.04%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.86'
$ws.Range('E40').Value = '  -2.73%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '308.46'
$ws.Range('E41').Value = '  -6.96%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '36.78'
$ws.Range('E42').Value = '  -1.95%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.66'
$ws.Range('E43').Value = '  -7.09%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.801'
$ws.Range('E44').Value = '  -13.06%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.594'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.77'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.52'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0923'
$ws.Range('E49').Value = '  -4.51%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.51'
$ws.Range('E50').Value = '  -5.55%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0514'
$ws.Range('E51').Value = '  -6.09%  '
